$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Direct value assignments (non-numeric-looking text, safe to assign directly)
$ws.Range('D2').Value = '42.939.99'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.361.55'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('E10').Value = '  -1.13%  '
$ws.Range('E11').Value = '  +3.28%  '
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('E13').Value = '  -3.28%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.726.93'
$ws.Range('E15').Value = '  +1.31%  '
$ws.Range('D16').Value = '2.354.98'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '42.891.60'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('E19').Value = '  -2.95%  '
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('E24').Value = '  -4.42%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('E28').Value = '  +15.32%  '
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  -11.93%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  +2.96%  '
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('E39').Value = '  +2.94%  '
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('E41').Value = '  -1.11%  '
$ws.Range('E42').Value = '  -5.04%  '
$ws.Range('D43').Value = '1.928.55'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('E45').Value = '  +3.28%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E46').Value = '  -9.26%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.588.21'
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E51').Value = '  +0.53%  '

# Numeric-looking strings in the Price column must be forced to remain text
Set-TextCell 'D4' '0.999'
Set-TextCell 'D5' '302.85'
Set-TextCell 'D6' '95.51'
Set-TextCell 'D9' '0.484'
Set-TextCell 'D10' '34.08'
Set-TextCell 'D12' '0.0784'
Set-TextCell 'D14' '6.71'
Set-TextCell 'D17' '0.790'
Set-TextCell 'D19' '11.93'
Set-TextCell 'D22' '67.98'
Set-TextCell 'D23' '235.01'
Set-TextCell 'D27' '24.47'
Set-TextCell 'D29' '9.31'
Set-TextCell 'D30' '31.99'
Set-TextCell 'D33' '17.46'
Set-TextCell 'D35' '128.07'
Set-TextCell 'D36' '1.84'
Set-TextCell 'D37' '0.104'
Set-TextCell 'D38' '4.30'
Set-TextCell 'D39' '2.83'
Set-TextCell 'D40' '2.26'
Set-TextCell 'D42' '20.97'
Set-TextCell 'D44' '0.0278'
Set-TextCell 'D45' '2.14'
Set-TextCell 'D46' '9.21'
Set-TextCell 'D47' '2.70'
Set-TextCell 'D49' '1.50'
Set-TextCell 'D50' '71.49'
Set-TextCell 'D51' '1.14'
